$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated shift start/end times and duration codes per row
$ws.Cells.Item(1, 3).Value = "07:00"
$ws.Cells.Item(1, 4).Value = "12:00"
$ws.Cells.Item(2, 3).Value = "15:00"
$ws.Cells.Item(2, 4).Value = "20:00"
$ws.Cells.Item(3, 3).Value = "07:00"
$ws.Cells.Item(3, 4).Value = "10:00"
$ws.Cells.Item(4, 3).Value = "12:00"
$ws.Cells.Item(4, 4).Value = "15:00"
$ws.Cells.Item(4, 5).Value = 60
$ws.Cells.Item(5, 3).Value = "07:00"
$ws.Cells.Item(5, 4).Value = "12:00"
$ws.Cells.Item(5, 5).Value = 80
$ws.Cells.Item(6, 3).Value = "12:00"
$ws.Cells.Item(6, 4).Value = "15:00"
$ws.Cells.Item(7, 3).Value = "14:00"
$ws.Cells.Item(7, 4).Value = "17:00"
$ws.Cells.Item(8, 4).Value = "10:00"
$ws.Cells.Item(8, 5).Value = 60
$ws.Cells.Item(9, 3).Value = "10:00"
$ws.Cells.Item(9, 4).Value = "15:00"
$ws.Cells.Item(9, 5).Value = 80
$ws.Cells.Item(10, 4).Value = "13:00"
$ws.Cells.Item(10, 5).Value = 80
$ws.Cells.Item(12, 3).Value = "07:00"
$ws.Cells.Item(12, 4).Value = "10:00"
$ws.Cells.Item(12, 5).Value = 60
$ws.Cells.Item(14, 3).Value = "07:00"
$ws.Cells.Item(14, 4).Value = "10:00"
$ws.Cells.Item(16, 3).Value = "07:00"
$ws.Cells.Item(16, 4).Value = "10:00"
$ws.Cells.Item(17, 3).Value = "10:00"
$ws.Cells.Item(17, 4).Value = "13:00"
$ws.Cells.Item(17, 5).Value = 60
$ws.Cells.Item(18, 3).Value = "13:00"
$ws.Cells.Item(18, 4).Value = "16:00"
$ws.Cells.Item(18, 5).Value = 60
$ws.Cells.Item(19, 3).Value = "08:00"
$ws.Cells.Item(19, 4).Value = "11:00"
$ws.Cells.Item(19, 5).Value = 60
$ws.Cells.Item(20, 3).Value = "13:00"
$ws.Cells.Item(20, 4).Value = "16:00"
$ws.Cells.Item(21, 3).Value = "15:00"
$ws.Cells.Item(21, 4).Value = "18:00"
$ws.Cells.Item(23, 3).Value = "07:00"
$ws.Cells.Item(23, 4).Value = "10:00"
$ws.Cells.Item(24, 3).Value = "11:00"
$ws.Cells.Item(24, 4).Value = "14:00"
$ws.Cells.Item(24, 5).Value = 60
$ws.Cells.Item(25, 3).Value = "11:00"
$ws.Cells.Item(25, 4).Value = "16:00"
$ws.Cells.Item(26, 3).Value = "14:00"
$ws.Cells.Item(26, 4).Value = "19:00"
$ws.Cells.Item(26, 5).Value = 80
$ws.Cells.Item(27, 3).Value = "08:00"
$ws.Cells.Item(27, 4).Value = "13:00"
$ws.Cells.Item(28, 3).Value = "12:00"
$ws.Cells.Item(28, 4).Value = "15:00"
$ws.Cells.Item(30, 3).Value = "15:00"
$ws.Cells.Item(30, 4).Value = "20:00"
$ws.Cells.Item(31, 4).Value = "15:00"
$ws.Cells.Item(31, 5).Value = 80
$ws.Cells.Item(32, 3).Value = "07:00"
$ws.Cells.Item(32, 4).Value = "10:00"
$ws.Cells.Item(32, 5).Value = 60
$ws.Cells.Item(33, 3).Value = "09:00"
$ws.Cells.Item(33, 4).Value = "12:00"
$ws.Cells.Item(34, 3).Value = "13:00"
$ws.Cells.Item(34, 4).Value = "18:00"
$ws.Cells.Item(35, 3).Value = "09:00"
$ws.Cells.Item(35, 4).Value = "14:00"
$ws.Cells.Item(36, 3).Value = "13:00"
$ws.Cells.Item(36, 5).Value = 80
$ws.Cells.Item(37, 3).Value = "08:00"
$ws.Cells.Item(37, 4).Value = "11:00"
$ws.Cells.Item(37, 5).Value = 60
$ws.Cells.Item(38, 3).Value = "11:00"
$ws.Cells.Item(38, 4).Value = "14:00"
$ws.Cells.Item(39, 3).Value = "14:00"
$ws.Cells.Item(39, 4).Value = "19:00"
$ws.Cells.Item(39, 5).Value = 80
$ws.Cells.Item(40, 3).Value = "14:00"
$ws.Cells.Item(40, 4).Value = "17:00"
$ws.Cells.Item(40, 5).Value = 60
$ws.Cells.Item(41, 3).Value = "12:00"
$ws.Cells.Item(41, 4).Value = "15:00"
$ws.Cells.Item(42, 3).Value = "12:00"
$ws.Cells.Item(42, 4).Value = "17:00"
$ws.Cells.Item(42, 5).Value = 80
$ws.Cells.Item(43, 3).Value = "15:00"
$ws.Cells.Item(43, 4).Value = "18:00"
$ws.Cells.Item(43, 5).Value = 60
$ws.Cells.Item(44, 3).Value = "14:00"
$ws.Cells.Item(44, 4).Value = "17:00"
$ws.Cells.Item(45, 3).Value = "07:00"
$ws.Cells.Item(45, 4).Value = "10:00"
$ws.Cells.Item(46, 3).Value = "14:00"
$ws.Cells.Item(46, 4).Value = "17:00"
$ws.Cells.Item(46, 5).Value = 60
$ws.Cells.Item(47, 3).Value = "07:00"
$ws.Cells.Item(47, 4).Value = "10:00"
$ws.Cells.Item(48, 3).Value = "13:00"
$ws.Cells.Item(48, 4).Value = "18:00"
$ws.Cells.Item(49, 3).Value = "15:00"
$ws.Cells.Item(49, 4).Value = "18:00"
$ws.Cells.Item(50, 3).Value = "13:00"
$ws.Cells.Item(50, 4).Value = "18:00"
$ws.Cells.Item(50, 5).Value = 80
$ws.Cells.Item(51, 3).Value = "14:00"
$ws.Cells.Item(51, 4).Value = "17:00"
$ws.Cells.Item(51, 5).Value = 60
$ws.Cells.Item(52, 3).Value = "15:00"
$ws.Cells.Item(52, 4).Value = "20:00"
$ws.Cells.Item(52, 5).Value = 80
$ws.Cells.Item(53, 3).Value = "12:00"
$ws.Cells.Item(53, 4).Value = "15:00"
$ws.Cells.Item(53, 5).Value = 60
$ws.Cells.Item(54, 3).Value = "08:00"
$ws.Cells.Item(54, 4).Value = "13:00"
$ws.Cells.Item(54, 5).Value = 80
$ws.Cells.Item(55, 3).Value = "07:00"
$ws.Cells.Item(55, 4).Value = "10:00"
$ws.Cells.Item(55, 5).Value = 60
$ws.Cells.Item(56, 3).Value = "15:00"
$ws.Cells.Item(56, 4).Value = "18:00"
$ws.Cells.Item(56, 5).Value = 60
$ws.Cells.Item(57, 3).Value = "15:00"
$ws.Cells.Item(57, 4).Value = "20:00"
$ws.Cells.Item(58, 3).Value = "14:00"
$ws.Cells.Item(58, 4).Value = "17:00"
$ws.Cells.Item(59, 3).Value = "08:00"
$ws.Cells.Item(59, 4).Value = "13:00"
$ws.Cells.Item(60, 3).Value = "09:00"
$ws.Cells.Item(60, 4).Value = "14:00"
$ws.Cells.Item(60, 5).Value = 80
$ws.Cells.Item(61, 3).Value = "08:00"
$ws.Cells.Item(61, 4).Value = "13:00"
$ws.Cells.Item(61, 5).Value = 80
$ws.Cells.Item(62, 3).Value = "09:00"
$ws.Cells.Item(62, 4).Value = "14:00"
$ws.Cells.Item(64, 3).Value = "14:00"
$ws.Cells.Item(64, 4).Value = "19:00"
$ws.Cells.Item(65, 3).Value = "15:00"
$ws.Cells.Item(65, 4).Value = "20:00"
$ws.Cells.Item(66, 3).Value = "11:00"
$ws.Cells.Item(66, 4).Value = "14:00"
$ws.Cells.Item(67, 3).Value = "12:00"
$ws.Cells.Item(67, 4).Value = "15:00"
$ws.Cells.Item(67, 5).Value = 60
$ws.Cells.Item(68, 3).Value = "11:00"
$ws.Cells.Item(68, 4).Value = "14:00"
$ws.Cells.Item(68, 5).Value = 60
$ws.Cells.Item(69, 3).Value = "12:00"
$ws.Cells.Item(69, 4).Value = "17:00"
$ws.Cells.Item(69, 5).Value = 80
$ws.Cells.Item(70, 3).Value = "14:00"
$ws.Cells.Item(70, 4).Value = "19:00"
$ws.Cells.Item(70, 5).Value = 80
$ws.Cells.Item(71, 3).Value = "11:00"
$ws.Cells.Item(71, 4).Value = "16:00"
$ws.Cells.Item(72, 4).Value = "15:00"
$ws.Cells.Item(72, 5).Value = 80
$ws.Cells.Item(73, 3).Value = "15:00"
$ws.Cells.Item(73, 4).Value = "20:00"
$ws.Cells.Item(73, 5).Value = 80
$ws.Cells.Item(74, 3).Value = "07:00"
$ws.Cells.Item(74, 4).Value = "12:00"
$ws.Cells.Item(75, 3).Value = "09:00"
$ws.Cells.Item(75, 4).Value = "14:00"
$ws.Cells.Item(75, 5).Value = 80
$ws.Cells.Item(76, 3).Value = "11:00"
$ws.Cells.Item(76, 4).Value = "14:00"
$ws.Cells.Item(76, 5).Value = 60
$ws.Cells.Item(77, 3).Value = "11:00"
$ws.Cells.Item(77, 4).Value = "14:00"
$ws.Cells.Item(78, 3).Value = "15:00"
$ws.Cells.Item(78, 4).Value = "20:00"
$ws.Cells.Item(79, 3).Value = "11:00"
$ws.Cells.Item(79, 4).Value = "16:00"
$ws.Cells.Item(80, 3).Value = "11:00"
$ws.Cells.Item(80, 4).Value = "14:00"
$ws.Cells.Item(81, 3).Value = "08:00"
$ws.Cells.Item(81, 4).Value = "11:00"
$ws.Cells.Item(81, 5).Value = 60
$ws.Cells.Item(82, 3).Value = "12:00"
$ws.Cells.Item(82, 4).Value = "17:00"
$ws.Cells.Item(83, 3).Value = "14:00"
$ws.Cells.Item(83, 4).Value = "19:00"
$ws.Cells.Item(83, 5).Value = 80
$ws.Cells.Item(85, 3).Value = "07:00"
$ws.Cells.Item(85, 4).Value = "10:00"
$ws.Cells.Item(86, 3).Value = "10:00"
$ws.Cells.Item(86, 4).Value = "15:00"
$ws.Cells.Item(87, 3).Value = "09:00"
$ws.Cells.Item(87, 4).Value = "12:00"
$ws.Cells.Item(88, 3).Value = "09:00"
$ws.Cells.Item(88, 4).Value = "14:00"
$ws.Cells.Item(88, 5).Value = 80
$ws.Cells.Item(90, 3).Value = "07:00"
$ws.Cells.Item(90, 4).Value = "12:00"
$ws.Cells.Item(91, 3).Value = "10:00"
$ws.Cells.Item(91, 4).Value = "13:00"
$ws.Cells.Item(92, 3).Value = "07:00"
$ws.Cells.Item(92, 4).Value = "12:00"
$ws.Cells.Item(93, 3).Value = "12:00"
$ws.Cells.Item(93, 4).Value = "15:00"
$ws.Cells.Item(94, 3).Value = "15:00"
$ws.Cells.Item(94, 4).Value = "20:00"
$ws.Cells.Item(95, 3).Value = "07:00"
$ws.Cells.Item(95, 4).Value = "12:00"
$ws.Cells.Item(96, 3).Value = "10:00"
$ws.Cells.Item(96, 4).Value = "15:00"
